{"js": "const pairs = [\n  [\"61\u00d734=2074\", \"84\u00d742=3528\"],\n  [\"19\u00d745=855\", \"61\u00d777=4697\"],\n  [\"78\u00d778=6084\", \"95\u00d730=2850\"],\n  [\"62\u00d713=806\", \"40\u00d767=2680\"],\n  [\"92\u00d790=8280\", \"76\u00d778=5928\"],\n  [\"83\u00d772=5976\", \"34\u00d757=1938\"],\n  [\"95\u00d740=3800\", \"33\u00d799=3267\"],\n  [\"96\u00d721=2016\", \"72\u00d796=6912\"],\n  [\"40\u00d731=1240\", \"22\u00d755=1210\"],\n  [\"58\u00d756=3248\", \"16\u00d736=576\"],\n  [\"42\u00d748=2016\", \"89\u00d754=4806\"],\n  [\"53\u00d732=1696\", \"66\u00d791=6006\"],\n  [\"81\u00d790=7290\", \"38\u00d758=2204\"],\n  [\"73\u00d785=6205\", \"60\u00d737=2220\"],\n  [\"18\u00d765=1170\", \"19\u00d790=1710\"],\n  [\"13\u00d771=923\", \"26\u00d787=2262\"],\n  [\"49\u00d785=4165\", \"57\u00d778=4446\"],\n  [\"48\u00d730=1440\", \"77\u00d754=4158\"],\n  [\"58\u00d746=2668\", \"17\u00d792=1564\"],\n  [\"75\u00d713=975\", \"32\u00d741=1312\"],\n  [\"83\u00d715=1245\", \"29\u00d796=2784\"],\n  [\"98\u00d787=8526\", \"15\u00d745=675\"],\n  [\"46\u00d749=2254\", \"17\u00d783=1411\"],\n  [\"68\u00d791=6188\", \"35\u00d774=2590\"],\n  [\"94\u00d720=1880\", \"53\u00d790=4770\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@('61\u00d734=2074', '84\u00d742=3528')\n    ,@('19\u00d745=855', '61\u00d777=4697')\n    ,@('78\u00d778=6084', '95\u00d730=2850')\n    ,@('62\u00d713=806', '40\u00d767=2680')\n    ,@('92\u00d790=8280', '76\u00d778=5928')\n    ,@('83\u00d772=5976', '34\u00d757=1938')\n    ,@('95\u00d740=3800', '33\u00d799=3267')\n    ,@('96\u00d721=2016', '72\u00d796=6912')\n    ,@('40\u00d731=1240', '22\u00d755=1210')\n    ,@('58\u00d756=3248', '16\u00d736=576')\n    ,@('42\u00d748=2016', '89\u00d754=4806')\n    ,@('53\u00d732=1696', '66\u00d791=6006')\n    ,@('81\u00d790=7290', '38\u00d758=2204')\n    ,@('73\u00d785=6205', '60\u00d737=2220')\n    ,@('18\u00d765=1170', '19\u00d790=1710')\n    ,@('13\u00d771=923', '26\u00d787=2262')\n    ,@('49\u00d785=4165', '57\u00d778=4446')\n    ,@('48\u00d730=1440', '77\u00d754=4158')\n    ,@('58\u00d746=2668', '17\u00d792=1564')\n    ,@('75\u00d713=975', '32\u00d741=1312')\n    ,@('83\u00d715=1245', '29\u00d796=2784')\n    ,@('98\u00d787=8526', '15\u00d745=675')\n    ,@('46\u00d749=2254', '17\u00d783=1411')\n    ,@('68\u00d791=6188', '35\u00d774=2590')\n    ,@('94\u00d720=1880', '53\u00d790=4770')\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $ok = $find.Execute(\n        $oldText,  # FindText\n        $false,    # MatchCase\n        $false,    # MatchWholeWord\n        $false,    # MatchWildcards\n        $false,    # MatchSoundsLike\n        $false,    # MatchAllWordForms\n        $true,     # Forward\n        1,         # Wrap (wdFindContinue)\n        $false,    # Format\n        $newText,  # ReplaceWith\n        2          # Replace (wdReplaceAll)\n    )\n    if (-not $ok) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
